$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 48686
$ws.Range("J17").Value = 50876.273
$ws.Range("L17").Value = 152628.819
$ws.Range("N17").Value = -152964.819
$ws.Range("H80").Value = 3893.353
$ws.Range("I80").Value = 468.06668
$ws.Range("J80").Value = 6597.5264
$ws.Range("K80").Value = 1404.20004
$ws.Range("L80").Value = 19792.5792
$ws.Range("M80").Value = -406.2000400000002
$ws.Range("N80").Value = -21788.5792
$ws.Range("H83").Value = 3893.353
$ws.Range("I83").Value = 468.06668
$ws.Range("J83").Value = 6597.5264
$ws.Range("K83").Value = 4212.60012
$ws.Range("L83").Value = 59377.7376
$ws.Range("M83").Value = 779.3998799999999
$ws.Range("N83").Value = -69361.73759999999
$ws.Range("H86").Value = 129000.375
$ws.Range("I86").Value = 336334.34
$ws.Range("J86").Value = 4600
$ws.Range("K86").Value = 336334.34
$ws.Range("L86").Value = 4600
$ws.Range("M86").Value = -335211.34
$ws.Range("N86").Value = -6846
$ws.Range("H89").Value = 129000.375
$ws.Range("I89").Value = 336334.34
$ws.Range("J89").Value = 4600
$ws.Range("K89").Value = 1681671.7
$ws.Range("L89").Value = 23000
$ws.Range("M89").Value = -1676055.7
$ws.Range("N89").Value = -34232
$ws.Range("H121").Value = 784.8570999999999
$ws.Range("I121").Value = 1331.3334
$ws.Range("J121").Value = 635.8182
$ws.Range("K121").Value = 3994.0002
$ws.Range("L121").Value = 1907.4546
$ws.Range("M121").Value = -2247.0002
$ws.Range("N121").Value = -5401.4546
$ws.Range("H125").Value = 14483.5
$ws.Range("I125").Value = 5399
$ws.Range("J125").Value = 17511.666
$ws.Range("K125").Value = 48591
$ws.Range("L125").Value = 157604.994
$ws.Range("M125").Value = -46131
$ws.Range("N125").Value = -162524.994
$ws.Range("H138").Value = 3682.0532
$ws.Range("I138").Value = 1614.5
$ws.Range("J138").Value = 4075.873
$ws.Range("K138").Value = 4843.5
$ws.Range("L138").Value = 12227.619
$ws.Range("M138").Value = 296.5
$ws.Range("N138").Value = -22507.619

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7864.0254
$ws.Range("I32").Value = 5899.0586
$ws.Range("K32").Value = 5899.0586
$ws.Range("M32").Value = -5612.0586
$ws.Range("H88").Value = 8077.375
$ws.Range("I88").Value = 14178
$ws.Range("J88").Value = 1976.75
$ws.Range("K88").Value = 14178
$ws.Range("L88").Value = 1976.75
$ws.Range("M88").Value = -13772
$ws.Range("N88").Value = -2788.75
$ws.Range("H91").Value = 8077.375
$ws.Range("I91").Value = 14178
$ws.Range("J91").Value = 1976.75
$ws.Range("K91").Value = 14178
$ws.Range("L91").Value = 1976.75
$ws.Range("M91").Value = -12774
$ws.Range("N91").Value = -4784.75
$ws.Range("H132").Value = 5398.757
$ws.Range("I132").Value = 1265.1904
$ws.Range("J132").Value = 10824.0625
$ws.Range("K132").Value = 3795.5712
$ws.Range("L132").Value = 32472.1875
$ws.Range("M132").Value = -1265.5712
$ws.Range("N132").Value = -37532.1875
$ws.Range("H140").Value = 60000
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 60000
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 60000
$ws.Range("M140").ClearContents()
$ws.Range("N140").Value = -70360

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2232.3513
$ws.Range("I86").Value = 2270.2942
$ws.Range("J86").Value = 1802.3334
$ws.Range("K86").Value = 2270.2942
$ws.Range("L86").Value = 1802.3334
$ws.Range("M86").Value = -1147.2942
$ws.Range("N86").Value = -4048.3334
$ws.Range("H89").Value = 2232.3513
$ws.Range("I89").Value = 2270.2942
$ws.Range("J89").Value = 1802.3334
$ws.Range("K89").Value = 11351.471
$ws.Range("L89").Value = 9011.666999999999
$ws.Range("M89").Value = -5735.471
$ws.Range("N89").Value = -20243.667

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H18").Value = 55470.4
$ws.Range("J18").Value = 55470.4
$ws.Range("L18").Value = 55470.4
$ws.Range("N18").Value = -55930.4
$ws.Range("H58").Value = 2167818.8
$ws.Range("I58").Value = 3499084
$ws.Range("J58").Value = 4512.4375
$ws.Range("K58").Value = 3499084
$ws.Range("L58").Value = 4512.4375
$ws.Range("M58").Value = -3498881
$ws.Range("N58").Value = -4918.4375
$ws.Range("H82").Value = 36181
$ws.Range("J82").Value = 36181
$ws.Range("L82").Value = 36181
$ws.Range("N82").Value = -36903
$ws.Range("H85").Value = 36181
$ws.Range("J85").Value = 36181
$ws.Range("L85").Value = 36181
$ws.Range("N85").Value = -38677
$ws.Range("H86").Value = 2714.6667
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 2714.6667
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 2714.6667
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = -4960.6667
$ws.Range("H89").Value = 2714.6667
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 2714.6667
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 13573.3335
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = -24805.3335
$ws.Range("H136").Value = 2167818.8
$ws.Range("I136").Value = 3499084
$ws.Range("J136").Value = 4512.4375
$ws.Range("K136").Value = 10497252
$ws.Range("L136").Value = 13537.3125
$ws.Range("M136").Value = -10494702
$ws.Range("N136").Value = -18637.3125

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 3017.2917
$ws.Range("I64").Value = 1273.7142
$ws.Range("J64").Value = 3735.2354
$ws.Range("K64").Value = 3821.1426
$ws.Range("L64").Value = 11205.7062
$ws.Range("M64").Value = -3551.1426
$ws.Range("N64").Value = -11745.7062
$ws.Range("H67").Value = 3017.2917
$ws.Range("I67").Value = 1273.7142
$ws.Range("J67").Value = 3735.2354
$ws.Range("K67").Value = 3821.1426
$ws.Range("L67").Value = 11205.7062
$ws.Range("M67").Value = -2885.1426
$ws.Range("N67").Value = -13077.7062
$ws.Range("H131").Value = 21912.156
$ws.Range("J131").Value = 33623.277
$ws.Range("L131").Value = 100869.831
$ws.Range("N131").Value = -110949.831
$ws.Range("H134").Value = 3666
$ws.Range("I134").Value = 3406.3572
$ws.Range("J134").Value = 4120.375
$ws.Range("K134").Value = 10219.0716
$ws.Range("L134").Value = 12361.125
$ws.Range("M134").Value = -5149.071599999999
$ws.Range("N134").Value = -22501.125

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 6626.5625
$ws.Range("I80").Value = 12151.8
$ws.Range("J80").Value = 4115.091
$ws.Range("K80").Value = 12151.8
$ws.Range("L80").Value = 4115.091
$ws.Range("M80").Value = -11153.8
$ws.Range("N80").Value = -6111.091
$ws.Range("H83").Value = 6626.5625
$ws.Range("I83").Value = 12151.8
$ws.Range("J83").Value = 4115.091
$ws.Range("K83").Value = 60759
$ws.Range("L83").Value = 20575.455
$ws.Range("M83").Value = -55767
$ws.Range("N83").Value = -30559.455
$ws.Range("H126").Value = 3002.5908
$ws.Range("I126").Value = 1912.5385
$ws.Range("K126").Value = 5737.6155
$ws.Range("M126").Value = -3267.6155
$ws.Range("H132").Value = 11844.833
$ws.Range("I132").Value = 17288.285
$ws.Range("J132").Value = 4224
$ws.Range("K132").Value = 51864.855
$ws.Range("L132").Value = 12672
$ws.Range("M132").Value = -49334.855
$ws.Range("N132").Value = -17732

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 932.75
$ws.Range("I16").Value = 932.75
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 932.75
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -762.75
$ws.Range("N16").ClearContents()
$ws.Range("H22").Value = 968.75
$ws.Range("H27").Value = 968.75
$ws.Range("H46").Value = 1080
$ws.Range("I46").Value = 900
$ws.Range("J46").Value = 1200
$ws.Range("K46").Value = 900
$ws.Range("L46").Value = 1200
$ws.Range("M46").Value = -712
$ws.Range("N46").Value = -1576
$ws.Range("H68").Value = 2071
$ws.Range("I68").Value = 1805.5714
$ws.Range("K68").Value = 1805.5714
$ws.Range("M68").Value = -1056.5714
$ws.Range("H71").Value = 2071
$ws.Range("I71").Value = 1805.5714
$ws.Range("K71").Value = 9027.857
$ws.Range("M71").Value = -5283.857
$ws.Range("H132").Value = 2895.0454
$ws.Range("I132").Value = 2155.923
$ws.Range("J132").Value = 3962.6667
$ws.Range("K132").Value = 6467.768999999999
$ws.Range("L132").Value = 11888.0001
$ws.Range("M132").Value = -3937.768999999999
$ws.Range("N132").Value = -16948.0001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3412.2354
$ws.Range("I62").Value = 3698
$ws.Range("J62").Value = 3293.1667
$ws.Range("K62").Value = 3698
$ws.Range("L62").Value = 3293.1667
$ws.Range("M62").Value = -3074
$ws.Range("N62").Value = -4541.1667
$ws.Range("H65").Value = 3412.2354
$ws.Range("I65").Value = 3698
$ws.Range("J65").Value = 3293.1667
$ws.Range("K65").Value = 18490
$ws.Range("L65").Value = 16465.8335
$ws.Range("M65").Value = -15370
$ws.Range("N65").Value = -22705.8335
$ws.Range("H81").Value = 3333.3333
$ws.Range("I81").Value = 2666.6667
$ws.Range("K81").Value = 5333.3334
$ws.Range("M81").Value = -4272.3334
$ws.Range("H84").Value = 3333.3333
$ws.Range("I84").Value = 2666.6667
$ws.Range("K84").Value = 26666.667
$ws.Range("M84").Value = -21362.667
$ws.Range("H132").Value = 4343.769
$ws.Range("I132").Value = 4156.9
$ws.Range("J132").Value = 4966.6665
$ws.Range("K132").Value = 12470.7
$ws.Range("L132").Value = 14899.9995
$ws.Range("M132").Value = -9940.699999999999
$ws.Range("N132").Value = -19959.9995
